# hopefully got rd2rus to run...
#
# Applies the "RUS_file format" table addition + Route_file format
# relabeling to the "endpoints" sheet, plus the K24:K28 unit-conversion
# helper column and the two highlight-color blocks.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("endpoints")

# --- text relabels -------------------------------------------------
$ws.Range("J4").Value = "location"
$ws.Range("F3").Value = "Route_file format"
$ws.Range("A9").Value  = "location"
$ws.Range("A10").Value = "keywords"

# --- new "RUS_file format" example table (N3:Q4) --------------------
$ws.Range("N3").Value = "RUS_file format"
$ws.Range("O4").Value = "user_ID"
$ws.Range("P4").Value = "stars"
$ws.Range("Q4").Value = "estimated stars"

# --- extend the "Route_file format" example row (F4:L4) -------------
$ws.Range("K4").Value = "number_pitches"
$ws.Range("L4").Value = "keywords"
$ws.Range("N4").Value = "route_ID"

# --- highlight blocks ------------------------------------------------
# green row under the "Route_file format" label (I3 keeps its existing
# quote-prefix style; adding the fill preserves that flag automatically)
$ws.Range("G3:K3").Interior.Color = 5287936   # 00B050 -> BGR 0x50B000
# orange row under the headline text
$ws.Range("G2:K2").Interior.Color = 49407     # FFC000 -> BGR 0x00C0FF

# --- rd2rus unit-conversion helper column (K24:K28) ------------------
$ws.Range("K24").Value = 100000
$ws.Range("K25").Value = 100
$ws.Range("K26").Formula = "=K24/K25"
$ws.Range("K27").Formula = "=K26*3"
$ws.Range("K28").Formula = "=K27/3600"

# widen the new helper column to fit "number_pitches"
$ws.Columns.Item(11).AutoFit() | Out-Null

# --- selection moved while exploring the new cells -------------------
$ws.Range("L32").Select()
